$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Update Si_No column (A12:A18) on Sheet1 - shift values, keep them as text
# (leading apostrophe forces text entry, matching the original quote-prefixed cells)
$ws1.Range("A12").Value = "'11"
$ws1.Range("A13").Value = "'12"
$ws1.Range("A14").Value = "'13"
$ws1.Range("A15").Value = "'14"
$ws1.Range("A16").Value = "'15"
$ws1.Range("A17").Value = "'16"
$ws1.Range("A18").Value = "'17"

# Update Sheet2 cell I2 value (plain number)
$ws2.Range("I2").Value = 7

# Update the selection on Sheet1 to D13
$ws1.Activate()
$ws1.Range("D13").Select()

# Update the selection on Sheet2 to I2
$ws2.Activate()
$ws2.Range("I2").Select()

# Re-activate Sheet1 since it's the tab that is selected in the target workbook
$ws1.Activate()
